$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 275.42307
$ws.Range("J9").Value = 284.5
$ws.Range("L9").Value = 284.5
$ws.Range("N9").Value = -622.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 11011.827
$ws.Range("J17").Value = 11438.63
$ws.Range("L17").Value = 34315.89
$ws.Range("N17").Value = -34651.89

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 71430104
$ws.Range("I92").Value = 71430104
$ws.Range("K92").Value = 71430104
$ws.Range("M92").Value = -71428856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1459.8889
$ws.Range("I101").Value = 1475
$ws.Range("K101").Value = 4425
$ws.Range("M101").Value = -2803

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 962.8946999999999
$ws.Range("J103").Value = 1070.4
$ws.Range("L103").Value = 3211.2
$ws.Range("N103").Value = -4383.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3413.99
$ws.Range("I138").Value = 1444.05
$ws.Range("J138").Value = 3906.475
$ws.Range("K138").Value = 4332.15
$ws.Range("L138").Value = 11719.425
$ws.Range("M138").Value = 807.8500000000004
$ws.Range("N138").Value = -21999.425

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 46063.25
$ws.Range("J24").Value = 46063.25
$ws.Range("L24").Value = 46063.25
$ws.Range("N24").Value = -46811.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4833.3335
$ws.Range("I63").Value = 4300
$ws.Range("J63").Value = 5500
$ws.Range("K63").Value = 4300
$ws.Range("L63").Value = 5500
$ws.Range("M63").Value = -3614
$ws.Range("N63").Value = -6872

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4833.3335
$ws.Range("I66").Value = 4300
$ws.Range("J66").Value = 5500
$ws.Range("K66").Value = 21500
$ws.Range("L66").Value = 27500
$ws.Range("M66").Value = -18068
$ws.Range("N66").Value = -34364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 46063.25
$ws.Range("J100").Value = 46063.25
$ws.Range("L100").Value = 46063.25
$ws.Range("N100").Value = -48227.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 58564.332
$ws.Range("J124").Value = 58564.332
$ws.Range("L124").Value = 58564.332
$ws.Range("N124").Value = -68384.33199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 100156.25
$ws.Range("J125").Value = 100156.25
$ws.Range("L125").Value = 100156.25
$ws.Range("N125").Value = -109996.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4007
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 4007
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 4007
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -6253

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4007
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 4007
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 20035
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -31267

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2117.5386
$ws.Range("I107").Value = 2029.4
$ws.Range("K107").Value = 2029.4
$ws.Range("M107").Value = -109.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4690000
$ws.Range("I6").Value = 4690000
$ws.Range("K6").Value = 4690000
$ws.Range("M6").Value = -4689887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 93.833336
$ws.Range("I7").Value = 40.77778
$ws.Range("K7").Value = 40.77778
$ws.Range("M7").Value = 72.22221999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 366834850
$ws.Range("I12").Value = 100000010
$ws.Range("J12").Value = 500252300
$ws.Range("K12").Value = 100000010
$ws.Range("L12").Value = 500252300
$ws.Range("M12").Value = -99999840
$ws.Range("N12").Value = -500252640

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 42963.34
$ws.Range("I86").Value = 37055.11
$ws.Range("J86").Value = 44625.03
$ws.Range("K86").Value = 37055.11
$ws.Range("L86").Value = 44625.03
$ws.Range("M86").Value = -35932.11
$ws.Range("N86").Value = -46871.03

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 42963.34
$ws.Range("I89").Value = 37055.11
$ws.Range("J89").Value = 44625.03
$ws.Range("K89").Value = 185275.55
$ws.Range("L89").Value = 223125.15
$ws.Range("M89").Value = -179659.55
$ws.Range("N89").Value = -234357.15

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1813.6
$ws.Range("I105").Value = 1626.2222
$ws.Range("K105").Value = 1626.2222
$ws.Range("M105").Value = 120.7778000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 86116.414
$ws.Range("I107").Value = 203119.8
$ws.Range("J107").Value = 2542.5715
$ws.Range("K107").Value = 203119.8
$ws.Range("L107").Value = 2542.5715
$ws.Range("M107").Value = -201199.8
$ws.Range("N107").Value = -6382.5715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1010.8333
$ws.Range("I136").Value = 739.0909
$ws.Range("K136").Value = 2217.2727
$ws.Range("M136").Value = 2882.7273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1891.25
$ws.Range("I140").Value = 1520
$ws.Range("K140").Value = 4560
$ws.Range("M140").Value = 620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 81614.664
$ws.Range("J120").Value = 81614.664
$ws.Range("L120").Value = 81614.664
$ws.Range("N120").Value = -91290.664

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4134.615
$ws.Range("I132").Value = 4104.2085
$ws.Range("K132").Value = 12312.6255
$ws.Range("M132").Value = -9782.625499999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 96616.836
$ws.Range("I132").Value = 104517.91
$ws.Range("K132").Value = 313553.73
$ws.Range("M132").Value = -311023.73

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1889.7037
$ws.Range("I136").Value = 1375.8667
$ws.Range("K136").Value = 4127.6001
$ws.Range("M136").Value = -1577.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4582.2666
$ws.Range("I62").Value = 3426
$ws.Range("J62").Value = 6316.6665
$ws.Range("K62").Value = 3426
$ws.Range("L62").Value = 6316.6665
$ws.Range("M62").Value = -2802
$ws.Range("N62").Value = -7564.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4582.2666
$ws.Range("I65").Value = 3426
$ws.Range("J65").Value = 6316.6665
$ws.Range("K65").Value = 17130
$ws.Range("L65").Value = 31583.3325
$ws.Range("M65").Value = -14010
$ws.Range("N65").Value = -37823.3325

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 60997.5
$ws.Range("J94").Value = 60997.5
$ws.Range("L94").Value = 60997.5
$ws.Range("N94").Value = -62799.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 11313.429
$ws.Range("J96").Value = 13050
$ws.Range("L96").Value = 13050
$ws.Range("N96").Value = -15796

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1775.1666
$ws.Range("I113").Value = 1530.2
$ws.Range("K113").Value = 4590.6
$ws.Range("M113").Value = -2420.6
